$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 153 ("「ブロック（90年代）」" entry) entirely; all rows below shift up by one.
$ws.Rows.Item(153).Delete()
